# Append new pupae-count rows (vial/treatment/time_hours=290/pupae) to Sheet1,
# mirroring the new experimental observations added at time_hours = 290.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for the newly recorded vials (vial, treatment, time_hours, pupae)
$newData = @(
    @(1, "conditioned", 290, 11),
    @(1, "unconditioned", 290, 19),
    @(2, "conditioned", 290, 12),
    @(2, "unconditioned", 290, 14),
    @(3, "conditioned", 290, 8),
    @(3, "unconditioned", 290, 4),
    @(4, "conditioned", 290, 14),
    @(4, "unconditioned", 290, 4),
    @(5, "conditioned", 290, 9),
    @(5, "unconditioned", 290, 7),
    @(6, "conditioned", 290, 5),
    @(6, "unconditioned", 290, 13),
    @(7, "conditioned", 290, 9),
    @(7, "unconditioned", 290, 21),
    @(8, "conditioned", 290, 8),
    @(8, "unconditioned", 290, 4),
    @(9, "conditioned", 290, 8),
    @(9, "unconditioned", 290, 4),
    @(10, "conditioned", 290, 6),
    @(10, "unconditioned", 290, 4),
    @(11, "conditioned", 290, 11),
    @(11, "unconditioned", 290, 12),
    @(12, "conditioned", 290, 21),
    @(12, "unconditioned", 290, 2),
    @(13, "conditioned", 290, 11),
    @(13, "unconditioned", 290, 5),
    @(15, "conditioned", 290, 5),
    @(14, "unconditioned", 290, 4),
    @(15, "conditioned", 290, 12),
    @(15, "unconditioned", 290, 4),
    @(16, "conditioned", 290, 4),
    @(16, "unconditioned", 290, 13),
    @(17, "conditioned", 290, 12),
    @(17, "unconditioned", 290, 2),
    @(18, "conditioned", 290, 13),
    @(18, "unconditioned", 290, 0),
    @(19, "conditioned", 290, 4),
    @(19, "unconditioned", 290, 5),
    @(20, "conditioned", 290, 13),
    @(20, "unconditioned", 290, 9),
    @(21, "conditioned", 290, 21),
    @(21, "unconditioned", 290, 5),
    @(22, "conditioned", 290, 10),
    @(22, "unconditioned", 290, 10),
    @(23, "conditioned", 290, 9),
    @(23, "unconditioned", 290, 3),
    @(24, "conditioned", 290, 30),
    @(24, "unconditioned", 290, 6),
    @(25, "conditioned", 290, 8),
    @(25, "unconditioned", 290, 0),
    @(26, "conditioned", 290, 8),
    @(26, "unconditioned", 290, 3),
    @(27, "conditioned", 290, 18),
    @(27, "unconditioned", 290, 4),
    @(28, "conditioned", 290, 24),
    @(28, "unconditioned", 290, 9),
    @(29, "conditioned", 290, 9),
    @(29, "unconditioned", 290, 7),
    @(30, "conditioned", 290, 8),
    @(30, "unconditioned", 290, 10),
    @(31, "conditioned", 290, 12),
    @(31, "unconditioned", 290, 12),
    @(32, "conditioned", 290, 3),
    @(32, "unconditioned", 290, 20),
    @(33, "conditioned", 290, 11),
    @(33, "unconditioned", 290, 16),
    @(34, "conditioned", 290, 11),
    @(34, "unconditioned", 290, 16)
)

# Find the first empty row below the existing data (existing data ends at row 1021)
$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $rec = $newData[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
}

$lastRow = $startRow + $newData.Count - 1

# Update the active selection/view to reflect where the new data was entered
$ws.Application.Goto($ws.Cells.Item($lastRow, 1), $true)
$excel.ActiveWindow.Zoom = 172
$ws.Range("D" + ($lastRow + 3)).Select()
